$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the newly added relationship opinions in column D (opinions of "Raiders")
$ws.Range("D2").Value = "For once does not trust someone at the door. Questions their authority but is quiet/shy/scared."
$ws.Range("D3").Value = "Readys for a fight. Thinks there are no possible good intentions."
$ws.Range("D4").Value = "Similar to Bob. Adrenaline starts pumping and she gets ready to defend."
$ws.Range("D5").Value = "Hateful/angry. These kinds of people woulded Sal."
$ws.Range("D6").Value = "Angry/nervous of harm coming to Hal. Tries to remain stoic and calm."

# Move the active selection from D3 to E2 (design doc moved to project)
$ws.Range("E2").Select()
